# Split the single run-on Bibliografia paragraph into separate sentences
# joined by manual line breaks (<w:br/>), one per reference entry.

$d = $word.ActiveDocument

$findText = "DENNIS, P. (2009). Produção lean simplificada. Bookman Editora." +
    "LEAN ENTERPRISE INSTITUTE (2007). Léxico lean – glossário ilustrado para praticantes do pensamento lean. v.2.0. São Paulo: Lean Institute Brasil." +
    "LIKER, J. K. (2005). O Modelo Toyota: 14 princípios de gestão do maior fabricante do mundo. Tradução  de Lene Belon Ribeiro. Porto Alegre: Bookman." +
    "ROSER, C. (2022). Tudo Sobre Produção Puxada: Projetando, Implementando e Mantendo Kanban, CONWIP e outros Sistemas Puxados na Produção Enxuta. AllAboutLean.com Publishing." +
    "ROTHER, M.; HARRIS, R. (2002). Criando fluxo contínuo. São Paulo, SP. Lean Institute Brasil." +
    "SHOOK, John; ROTHER, Mike. Manual. Aprendendo a enxergar. Leam Institute Brasil. São Paulo: IMAM, s/d." +
    "TUBINO, D. F. (2015). Manufatura enxuta como estratégia de produção. Editora Atlas SA." +
    "WOMACK, James P.; JONES, Daniel T. A Mentalidade enxuta nas empresas. Rio de Janeiro: Campus, 1998. " +
    "WOMACK, James P.; JONES, Daniel T. Lean Thinking: Banish Waste and Create Wealth in Your Corporation. Free Press, 2010."

$replaceText = "DENNIS, P. (2009). Produção lean simplificada. Bookman Editora." + "^l" +
    "LEAN ENTERPRISE INSTITUTE (2007). Léxico lean – glossário ilustrado para praticantes do pensamento lean. v.2.0. São Paulo: Lean Institute Brasil." + "^l" +
    "LIKER, J. K. (2005). O Modelo Toyota: 14 princípios de gestão do maior fabricante do mundo. Tradução  de Lene Belon Ribeiro. Porto Alegre: Bookman." + "^l" +
    "ROSER, C. (2022). Tudo Sobre Produção Puxada: Projetando, Implementando e Mantendo Kanban, CONWIP e outros Sistemas Puxados na Produção Enxuta. AllAboutLean.com Publishing." + "^l" +
    "ROTHER, M.; HARRIS, R. (2002). Criando fluxo contínuo. São Paulo, SP. Lean Institute Brasil." + "^l" +
    "SHOOK, John; ROTHER, Mike. Manual. Aprendendo a enxergar. Leam Institute Brasil. São Paulo: IMAM, s/d." + "^l" +
    "TUBINO, D. F. (2015). Manufatura enxuta como estratégia de produção. Editora Atlas SA." + "^l" +
    "WOMACK, James P.; JONES, Daniel T. A Mentalidade enxuta nas empresas. Rio de Janeiro: Campus, 1998. " + "^l" +
    "WOMACK, James P.; JONES, Daniel T. Lean Thinking: Banish Waste and Create Wealth in Your Corporation. Free Press, 2010."

$range = $d.Content
$found = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)

Write-Output "Replaced bibliography run-on text with line-broken entries: $found"
